# 13.1.1.xlsx — add the 2020 column (Q) with disaster-death data, mirroring
# the formatting of the existing 2019 column (P), move the active selection,
# and let the sheet's dimension / row spans update automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Row 3: bottom-border-only separator row, no value, just copy style ---
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial($xlPasteFormats)

# --- Data rows: (row, value, isDash) ------------------------------------
# isDash rows get the literal "-" placeholder (shared string), others get
# a numeric value.
$rows = @(
    @{R=4;  V=2020; Dash=$false},
    @{R=5;  V=51;   Dash=$false},
    @{R=6;  V=29;   Dash=$false},
    @{R=7;  V=22;   Dash=$false},
    @{R=8;  V=5;    Dash=$false},
    @{R=9;  V=3;    Dash=$false},
    @{R=10; V=2;    Dash=$false},
    @{R=11; V=15;   Dash=$false},
    @{R=12; V=9;    Dash=$false},
    @{R=13; V=5;    Dash=$false},
    @{R=14; V=0;    Dash=$true},
    @{R=15; V=0;    Dash=$true},
    @{R=16; V=0;    Dash=$true},
    @{R=17; V=0;    Dash=$true},
    @{R=18; V=0;    Dash=$true},
    @{R=19; V=0;    Dash=$true},
    @{R=20; V=7;    Dash=$false},
    @{R=21; V=7;    Dash=$false},
    @{R=22; V=0;    Dash=$true},
    @{R=23; V=0;    Dash=$true},
    @{R=24; V=0;    Dash=$true},
    @{R=25; V=0;    Dash=$true},
    @{R=26; V=24;   Dash=$false},
    @{R=27; V=10;   Dash=$false},
    @{R=28; V=14;   Dash=$false},
    @{R=29; V=0;    Dash=$true},
    @{R=30; V=0;    Dash=$true},
    @{R=31; V=0;    Dash=$true},
    @{R=32; V=0;    Dash=$true},
    @{R=33; V=0;    Dash=$true},
    @{R=34; V=0;    Dash=$true}
)

foreach ($row in $rows) {
    $r = $row.R
    $srcCell = $ws.Range("P" + $r)
    $dstCell = $ws.Range("Q" + $r)

    # Copy P<r>'s number format / font / alignment / borders onto Q<r>.
    $srcCell.Copy()
    $dstCell.PasteSpecial($xlPasteFormats)

    if ($row.Dash) {
        $dstCell.Value = "-"
    } else {
        $dstCell.Value = $row.V
    }
}

# --- Move the active selection to K18, as in the authored workbook ------
$ws.Range("K18").Select()

$wb.Save()
